$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: "Test Exp 10" -> "Test Exp 11"
$ws.Range("B12").Value = "Test Exp 11"

# New row 13, mirroring row 12's former content but with "Test Exp 12"
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "Test Exp 12"
$ws.Range("C13").Value = 30
$ws.Range("D13").Value = 0.25
$ws.Range("E13").Value = "Local"
$ws.Range("F13").Value = -1
$ws.Range("G13").Value = "28*28"
$ws.Range("H13").Value = "32*32"
$ws.Range("I13").Value = "3,4,5"

# Match style (left-aligned) used by the rest of the table rows (I column
# intentionally left at default style, matching every other row's I cell)
$ws.Range("A13:H13").HorizontalAlignment = -4131

# Update selection to mirror the saved workbook state
$ws.Range("D14").Select()
